$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the requested frequency (B7) to reflect the restored design intent
$ws.Range("B7").Value = 51.5

# Force recalculation so dependent formulas (B9, B10) refresh
$excel.Calculate()

# Restore the selection/active cell as captured in the saved view state
$ws.Range("F8").Select()
